# Generate Report for Handback
# Updates the localization-status workbook to reflect that both files have
# been handed back and are in sync with en-US: fills in the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns on the
# zh-cn and de-de sheets, flips the Overview Status to the handed-back
# message, and widens the columns that now hold longer content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("J2").Value = "624f1614-dda0-476e-b4cb-c830f7275612.2421fe3189219b8aa933c0162d781b99e14d28ab.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 06:29:13"

$wsZh.Range("J3").Value = "dbb31382-a233-4f90-9595-e2cc6c104265.5f950d577d594652dc38748aeb572786b537e9fe.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 06:29:13"

# Rebuild the hyperlinks so the new "Latest Target File" cells (I2/I3) link
# back to the source .md files, same as column A already does.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/624f1614-dda0-476e-b4cb-c830f7275612.md", "", "", "624f1614-dda0-476e-b4cb-c830f7275612.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/624f1614-dda0-476e-b4cb-c830f7275612.md", "", "", "624f1614-dda0-476e-b4cb-c830f7275612.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/dbb31382-a233-4f90-9595-e2cc6c104265.md", "", "", "dbb31382-a233-4f90-9595-e2cc6c104265.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/dbb31382-a233-4f90-9595-e2cc6c104265.md", "", "", "dbb31382-a233-4f90-9595-e2cc6c104265.md")

# Widen columns to fit the newly-populated data.
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J2").Value = "624f1614-dda0-476e-b4cb-c830f7275612.2421fe3189219b8aa933c0162d781b99e14d28ab.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 06:29:20"

$wsDe.Range("J3").Value = "dbb31382-a233-4f90-9595-e2cc6c104265.5f950d577d594652dc38748aeb572786b537e9fe.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 06:29:20"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/624f1614-dda0-476e-b4cb-c830f7275612.md", "", "", "624f1614-dda0-476e-b4cb-c830f7275612.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/624f1614-dda0-476e-b4cb-c830f7275612.md", "", "", "624f1614-dda0-476e-b4cb-c830f7275612.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/dbb31382-a233-4f90-9595-e2cc6c104265.md", "", "", "dbb31382-a233-4f90-9595-e2cc6c104265.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119778ffc16dabaad412b8d38368102f0518ddde/e2e/dbb31382-a233-4f90-9595-e2cc6c104265.md", "", "", "dbb31382-a233-4f90-9595-e2cc6c104265.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------------
# Overview sheet - widen the per-language status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

Write-Host "Handback report generated."
